$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "Save" in H1, matching the style of G1 (bold/centered/bordered header style)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add values in H2 and H3
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
